# Update Repositories.xlsx to reflect recent updates:
# Insert two new repository entries (NFDI-MatWerk Repository and the
# Inorganic Crystal Structure Database) at the top of the "Materials
# Science" block on the Repositories sheet, pushing the existing rows
# down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repositories")
$ws.Activate()

# Insert two blank rows above the current row 50 - everything that was
# row 50 onward (through row 85) shifts down to rows 52-87.
$ws.Rows("50:51").Insert()

# New row 50: KIT - NFDI-MatWerk Repository
$ws.Range("A50").Value2 = "KIT"
$ws.Range("B50").Value2 = "NFDI-MatWerk Repository"
$ws.Range("C50").Value2 = "Repository"
$ws.Range("D50").Value2 = "Materials Science"
$ws.Range("E50").Value2 = "yes"
$ws.Range("F50").Value2 = "yes"
$ws.Range("G50").Value2 = "free"
$ws.Range("I50").Value2 = "https://matwerk.datamanager.kit.edu"
$ws.Range("J50").Value2 = "Data and metadata repository for NFDI-MatWerk. File size limit for the metadata repository is 10MB."

# New row 51: FIZ Karlsruhe - Inorganic Crystal Structure Database
$ws.Range("A51").Value2 = "FIZ Karlsruhe"
$ws.Range("B51").Value2 = "Inorganic Crystal Structure Database"
$ws.Range("C51").Value2 = "Database"
$ws.Range("D51").Value2 = "Materials Science"
$ws.Range("E51").Value2 = "no"
$ws.Range("F51").Value2 = "yes"
$ws.Range("I51").Value2 = "https://icsd.products.fiz-karlsruhe.de"
$ws.Range("J51").Value2 = "Database for completely identified inorganic crystal structures. Access requires a license."

# Restore the view state roughly where the author left it: scrolled down
# to around row 21, with cell J51 (the new comment cell) selected.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J51").Select()
